# Applies numeric corrections to the per-job Leve profit tables (H:N columns)
# across all 8 sheets, per the scheduled-runner price refresh.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 1465.6428
$ws.Range("J53").Value = 4631.3335
$ws.Range("L53").Value = 4631.3335
$ws.Range("N53").Value = -5905.3335
$ws.Range("H98").Value = 1468.4
$ws.Range("I98").Value = 1629.8667
$ws.Range("K98").Value = 1629.8667
$ws.Range("M98").Value = -131.8667
$ws.Range("H113").Value = 5128
$ws.Range("I113").Value = 4780.4
$ws.Range("J113").Value = 5417.6665
$ws.Range("K113").Value = 4780.4
$ws.Range("L113").Value = 5417.6665
$ws.Range("M113").Value = -1526.4
$ws.Range("N113").Value = -11925.6665
$ws.Range("H122").Value = 1468.4
$ws.Range("I122").Value = 1629.8667
$ws.Range("K122").Value = 4889.6001
$ws.Range("M122").Value = -2439.6001
$ws.Range("H128").Value = 130000
$ws.Range("J128").Value = 130000
$ws.Range("L128").Value = 130000
$ws.Range("N128").Value = -139960
$ws.Range("H130").Value = 116778.5
$ws.Range("J130").Value = 116778.5
$ws.Range("L130").Value = 116778.5
$ws.Range("N130").Value = -126818.5
$ws.Range("H132").Value = 1278.4865
$ws.Range("I132").Value = 1316.9714
$ws.Range("J132").Value = 605
$ws.Range("K132").Value = 3950.9142
$ws.Range("L132").Value = 1815
$ws.Range("M132").Value = -1420.9142
$ws.Range("N132").Value = -6875
$ws.Range("H138").Value = 3454.35
$ws.Range("I138").Value = 2377.158
$ws.Range("J138").Value = 3707.0247
$ws.Range("K138").Value = 7131.474
$ws.Range("L138").Value = 11121.0741
$ws.Range("M138").Value = -1991.474
$ws.Range("N138").Value = -21401.0741
$ws.Range("H141").Value = 6465.4
$ws.Range("I141").Value = 6506.087
$ws.Range("K141").Value = 19518.261
$ws.Range("M141").Value = -14338.261

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15801726
$ws.Range("I32").Value = 17511728
$ws.Range("J32").Value = 5883719
$ws.Range("K32").Value = 17511728
$ws.Range("L32").Value = 5883719
$ws.Range("M32").Value = -17511441
$ws.Range("N32").Value = -5884293
$ws.Range("H61").Value = 4441
$ws.Range("I61").Value = 4527.2144
$ws.Range("J61").Value = 4199.6
$ws.Range("K61").Value = 4527.2144
$ws.Range("L61").Value = 4199.6
$ws.Range("M61").Value = -4315.2144
$ws.Range("N61").Value = -4623.6
$ws.Range("H102").Value = 2293.182
$ws.Range("I102").Value = 1858.3334
$ws.Range("K102").Value = 1858.3334
$ws.Range("M102").Value = -236.3334
$ws.Range("H110").Value = 1773.3334
$ws.Range("I110").Value = 785
$ws.Range("K110").Value = 785
$ws.Range("M110").Value = 1260
$ws.Range("H132").Value = 2497.5208
$ws.Range("I132").Value = 2287.6904
$ws.Range("K132").Value = 6863.0712
$ws.Range("M132").Value = -4333.0712
$ws.Range("H136").Value = 4441
$ws.Range("I136").Value = 4527.2144
$ws.Range("J136").Value = 4199.6
$ws.Range("K136").Value = 13581.6432
$ws.Range("L136").Value = 12598.8
$ws.Range("M136").Value = -11031.6432
$ws.Range("N136").Value = -17698.8

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3077.5557
$ws.Range("I20").Value = 3197.1428
$ws.Range("J20").Value = 2659
$ws.Range("K20").Value = 3197.1428
$ws.Range("L20").Value = 2659
$ws.Range("M20").Value = -2950.1428
$ws.Range("N20").Value = -3153
$ws.Range("H92").Value = 50401
$ws.Range("J92").Value = 50401
$ws.Range("L92").Value = 50401
$ws.Range("N92").Value = -55393
$ws.Range("H94").Value = 988.2105
$ws.Range("I94").Value = 868.1539
$ws.Range("K94").Value = 868.1539
$ws.Range("M94").Value = -417.1539
$ws.Range("H99").Value = 2745.2632
$ws.Range("I99").Value = 2510.8572
$ws.Range("K99").Value = 2510.8572
$ws.Range("M99").Value = -1012.8572
$ws.Range("H105").Value = 2730.9375
$ws.Range("I105").Value = 2406.7856
$ws.Range("K105").Value = 2406.7856
$ws.Range("M105").Value = -659.7856000000002
$ws.Range("H134").Value = 2566802.8
$ws.Range("I134").Value = 2780286.5
$ws.Range("K134").Value = 8340859.5
$ws.Range("M134").Value = -8338324.5

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1726.4839
$ws.Range("I31").Value = 1201.3334
$ws.Range("J31").Value = 2218.8125
$ws.Range("K31").Value = 1201.3334
$ws.Range("L31").Value = 2218.8125
$ws.Range("M31").Value = -906.3334
$ws.Range("N31").Value = -2808.8125
$ws.Range("H34").Value = 1726.4839
$ws.Range("I34").Value = 1201.3334
$ws.Range("J34").Value = 2218.8125
$ws.Range("K34").Value = 1201.3334
$ws.Range("L34").Value = 2218.8125
$ws.Range("M34").Value = -999.3334
$ws.Range("N34").Value = -2622.8125
$ws.Range("H95").Value = 64881
$ws.Range("J95").Value = 64881
$ws.Range("L95").Value = 64881
$ws.Range("N95").Value = -70373
$ws.Range("H107").Value = 614.25
$ws.Range("I107").Value = 419.3
$ws.Range("J107").Value = 1589
$ws.Range("K107").Value = 419.3
$ws.Range("L107").Value = 1589
$ws.Range("M107").Value = 1500.7
$ws.Range("N107").Value = -5429
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H132").Value = 3587.5454
$ws.Range("I132").Value = 3164.6316
$ws.Range("K132").Value = 9493.8948
$ws.Range("M132").Value = -6963.8948
$ws.Range("H134").Value = 3780.3333
$ws.Range("I134").Value = 3812.2856
$ws.Range("J134").Value = 3333
$ws.Range("K134").Value = 11436.8568
$ws.Range("L134").Value = 9999
$ws.Range("M134").Value = -8901.856800000001
$ws.Range("N134").Value = -15069

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 399.625
$ws.Range("I5").Value = 406.7143
$ws.Range("J5").Value = 350
$ws.Range("K5").Value = 1220.1429
$ws.Range("L5").Value = 1050
$ws.Range("M5").Value = -1108.1429
$ws.Range("N5").Value = -1274
$ws.Range("H68").Value = 2151.7576
$ws.Range("J68").Value = 2324.6
$ws.Range("L68").Value = 6973.799999999999
$ws.Range("N68").Value = -8595.799999999999
$ws.Range("H71").Value = 2151.7576
$ws.Range("J71").Value = 2324.6
$ws.Range("L71").Value = 20921.4
$ws.Range("N71").Value = -29033.4
$ws.Range("H129").Value = 2187.4546
$ws.Range("J129").Value = 2187.4546
$ws.Range("L129").Value = 6562.3638
$ws.Range("N129").Value = -16562.3638
$ws.Range("H135").Value = 399.625
$ws.Range("I135").Value = 406.7143
$ws.Range("J135").Value = 350
$ws.Range("K135").Value = 3660.4287
$ws.Range("L135").Value = 3150
$ws.Range("M135").Value = -1125.4287
$ws.Range("N135").Value = -8220
$ws.Range("H140").Value = 1789.5555
$ws.Range("I140").Value = 1434.3334
$ws.Range("K140").Value = 4303.0002
$ws.Range("M140").Value = 876.9997999999996

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 22750
$ws.Range("J40").Value = 25500
$ws.Range("L40").Value = 25500
$ws.Range("N40").Value = -25802
$ws.Range("H97").Value = 1227.9286
$ws.Range("I97").Value = 1084.3334
$ws.Range("J97").Value = 1658.7142
$ws.Range("K97").Value = 1084.3334
$ws.Range("L97").Value = 1658.7142
$ws.Range("M97").Value = -588.3334
$ws.Range("N97").Value = -2650.7142
$ws.Range("H102").Value = 2267.375
$ws.Range("I102").Value = 2098.6924
$ws.Range("K102").Value = 2098.6924
$ws.Range("M102").Value = -476.6923999999999
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H132").Value = 3940.7222
$ws.Range("I132").Value = 3662.0833
$ws.Range("J132").Value = 4498
$ws.Range("K132").Value = 10986.2499
$ws.Range("L132").Value = 13494
$ws.Range("M132").Value = -8456.249899999999
$ws.Range("N132").Value = -18554
$ws.Range("H136").Value = 57296.637
$ws.Range("J136").Value = 57296.637
$ws.Range("L136").Value = 171889.911
$ws.Range("N136").Value = -176989.911

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6857.467
$ws.Range("I132").Value = 6921.3447
$ws.Range("K132").Value = 20764.0341
$ws.Range("M132").Value = -18234.0341
$ws.Range("H136").Value = 35272.637
$ws.Range("I136").Value = 38400
$ws.Range("J136").Value = 3999
$ws.Range("K136").Value = 115200
$ws.Range("L136").Value = 11997
$ws.Range("M136").Value = -112650
$ws.Range("N136").Value = -17097

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 109989.5
$ws.Range("J16").Value = 109989.5
$ws.Range("L16").Value = 109989.5
$ws.Range("N16").Value = -110573.5
$ws.Range("H132").Value = 1888.3208
$ws.Range("I132").Value = 1764.5238
$ws.Range("K132").Value = 5293.5714
$ws.Range("M132").Value = -2763.5714
$ws.Range("H136").Value = 34305.53
$ws.Range("I136").Value = 2820.3157
$ws.Range("J136").Value = 80322.38
$ws.Range("K136").Value = 8460.947100000001
$ws.Range("L136").Value = 240967.14
$ws.Range("M136").Value = -5910.947100000001
$ws.Range("N136").Value = -246067.14

